# Revert "Add snRNAseq-10xGenomics-v2 to the scrnaseq assays"
#
# The original commit added a "snRNAseq-10xGenomics-v2" entry to the
# "assay_type list" sheet (row 3, between scRNAseq-10xGenomics-v3 and
# scRNAseq) and widened the K2:K1048576 data validation on "Export as TSV"
# from 'assay_type list'!$A$1:$A$6 to $A$1:$A$7 to include it.
#
# Reverting: remove that row from the lookup list and shrink the data
# validation range back down. Excel takes care of compacting the shared
# string table (and therefore every other <v> index that shifts as a
# result) automatically when the file is saved.

$wb = $excel.ActiveWorkbook

# 1. Remove the "snRNAseq-10xGenomics-v2" row from the "assay_type list" sheet.
$assayTypeList = $wb.Worksheets.Item("assay_type list")
$assayTypeList.Range("A3").EntireRow.Delete()

# 2. Shrink the assay_type data validation range on the main sheet back to
#    $A$1:$A$6 (was $A$1:$A$7 while the extra entry existed).
$tsv = $wb.Worksheets.Item("Export as TSV")
$dv = $tsv.Range("K2:K1048576").Validation
$dv.Modify(3, 1, 1, "'assay_type list'!`$A`$1:`$A`$6")
